$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row appended by the 2025-10-04 run.
$row = 33

# Force the date column to be stored as literal text (matching the rest of
# the "Date" column, which holds plain MM/DD/YYYY strings rather than real
# dates) instead of letting Excel auto-parse "10/04/2025" into a date
# serial. Resetting the Style back to "Normal" afterwards keeps the cell
# free of any explicit formatting, just like its neighbours above.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "10/04/2025"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = 0.1465718972668775
$ws.Range("C$row").Value = 0.8534281027331225
